# Apply "added final result calcs after adjusting ML labels"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row labels: TSLA row is now the "Log. Reg." model, BTC row is now the "SVC" model ---
$ws.Range("A3").Value = "TSLA (Log. Reg.)"
$ws.Range("A4").Value = "BTC (SVC)"

# --- Updated gain figures (B:D) for each row; E column formulas (AVERAGE) recalc automatically ---
$ws.Range("B2").Value = -0.0037588512066153434
$ws.Range("C2").Value = -0.016233611113976476
$ws.Range("D2").Value = 0.01939031661716495

$ws.Range("B3").Value = 0.27427033219897146
$ws.Range("C3").Value = 0.2567049671383368
$ws.Range("D3").Value = 0.2647301032011935

$ws.Range("B4").Value = -0.11411730577160228
$ws.Range("C4").Value = -0.17394731033348698
$ws.Range("D4").Value = -0.09777181474141194

# --- Colour the figures: green for gains, red for losses (matches existing workbook convention) ---
$green = 5287936   # RGB(0,176,80)  -> FF00B050
$red   = 255       # RGB(255,0,0)   -> FFFF0000

$ws.Range("B2").Font.Color = $red
$ws.Range("C2").Font.Color = $red
$ws.Range("D2").Font.Color = $green
$ws.Range("E2").Font.Color = $red

$ws.Range("B3").Font.Color = $green
$ws.Range("C3").Font.Color = $green
$ws.Range("D3").Font.Color = $green
$ws.Range("E3").Font.Color = $green

$ws.Range("B4").Font.Color = $red
$ws.Range("C4").Font.Color = $red
$ws.Range("D4").Font.Color = $red
$ws.Range("E4").Font.Color = $red

$ws.Range("B5").Font.Color = $green
$ws.Range("C5").Font.Color = $green
$ws.Range("D5").Font.Color = $green

# --- Restore the cursor position left by the author ---
$ws.Range("C19").Select()
